$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove duplicated "MCT-1A-Gestão integrada" entries in column B (rows 2 and 3)
$ws.Range("B2").Value = "-"
$ws.Range("B3").Value = "-"

# Fix mistakenly duplicated "MEC-1B-Gestao Intregrada" entries in column E (rows 6 and 7)
# which should actually be "MCT-1A-Gestão integrada"
$ws.Range("E6").Value = "MCT-1A-Gestão integrada"
$ws.Range("E7").Value = "MCT-1A-Gestão integrada"

# Remove duplicated "MEC-1A-Gestao Integrada" entries in column B (rows 12 and 14)
$ws.Range("B12").Value = "-"
$ws.Range("B14").Value = "-"
